$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 4.621579
$ws.Range("N2").Value = 13.864737
$ws.Range("O2").Value = 0.1778708528171788
$ws.Range("P2").Value = 0.1778708528171788
$ws.Range("Q2").Value = 19.77176968569167
$ws.Range("R2").Value = 177.945927171225
$ws.Range("S2").Value = 0.1718806088494653
$ws.Range("T2").Value = 0.1718806088494653

# Row 3
$ws.Range("N3").Value = 46.543441
$ws.Range("O3").Value = 0.5971062807549863
$ws.Range("P3").Value = 0.5971062807549863
$ws.Range("Q3").Value = 66.37314475071389
$ws.Range("R3").Value = 597.358302756425
$ws.Range("S3").Value = 0.5769972396179723
$ws.Range("T3").Value = 0.5769972396179723

# Row 4
$ws.Range("O4").Value = 0.2250228664278349
$ws.Range("P4").Value = 0.2250228664278349
$ws.Range("S4").Value = 0.2174446609665816
$ws.Range("T4").Value = 0.2174446609665815

# Row 5
$ws.Range("M5").Value = 4.621579
$ws.Range("N5").Value = 13.864737
$ws.Range("O5").Value = 0.1778708528171788
$ws.Range("P5").Value = 0.1778708528171788
$ws.Range("Q5").Value = 0.6890697262683333
$ws.Range("R5").Value = 6.201627536415
$ws.Range("S5").Value = 0.00599024396771352
$ws.Range("T5").Value = 0.005990243967713519

# Row 6
$ws.Range("N6").Value = 46.543441
$ws.Range("O6").Value = 0.5971062807549863
$ws.Range("P6").Value = 0.5971062807549863
$ws.Range("S6").Value = 0.02010904113701401
$ws.Range("T6").Value = 0.020109041137014

# Row 7
$ws.Range("O7").Value = 0.2250228664278349
$ws.Range("P7").Value = 0.2250228664278349
$ws.Range("S7").Value = 0.007578205461253396
$ws.Range("T7").Value = 0.007578205461253394
